# Update order summary after successful item deletion
# Adds rows 46-51 to Sheet1 (Order_Items) with the new order item data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(45, 43, "PAMP001",  3, 22.5,               0),
    @(46, 44, "PAMP001", 18, 135,                0),
    @(47, 44, "PAMR2234",24, 382.7796610169491,  1),
    @(48, 45, "PMP234",   1, 2446.035874439462,  0),
    @(49, 45, "PAMP0000", 1, 1,                  1),
    @(50, 45, "CAND234",  1, 896.8609865470852,  0)
)

$startRow = 46
for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $startRow + $i
    $row = $data[$i]

    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
}
